$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Header row additions: Wins, Losses, Ties
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy header style from an existing header cell (AC1) to new header cells
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122) # xlPasteFormats

# Fill data rows 2-42 with Wins=96, Losses=66, Ties=0
for ($r = 2; $r -le 42; $r++) {
    $ws.Cells.Item($r, 30).Value = 96  # AD
    $ws.Cells.Item($r, 31).Value = 66  # AE
    $ws.Cells.Item($r, 32).Value = 0   # AF
}
